# Adds a new "#" / "{{index+1}}" column to the PlanData template.
#
# The sheet has three independent regions that each shift one column to the
# right starting at column B (a new column is inserted there), while the
# decorative A:G number grid in rows 6-11 is intentionally left untouched:
#   - row 1:  the merged title cell F1:J1  -> G1:K1
#   - rows 3-4: the template header/row   B:J -> C:K (new B gets "#" / "{{index+1}}")
#   - rows 5-9: the small helper matrix   (D:E on row5, H:J on rows 3-9) -> one col right
#
# Because this runtime's Range.Insert()/Columns.Insert() always shifts the
# *entire* column (which would also drag the unrelated A:G number grid
# along), the move is instead done cell-by-cell, right-to-left, using
# Copy+PasteSpecial(xlPasteFormats) to carry styles and explicit .Value
# writes to carry data - fully equivalent to the end result of the
# insert, without disturbing the untouched region.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Move-Cell($destAddr, $srcAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($destAddr).PasteSpecial($xlPasteFormats)
}

# ---- Row 1: F1:J1 -> G1:K1 (merged title cell) ----
Move-Cell "K1" "J1"
Move-Cell "J1" "I1"
Move-Cell "I1" "H1"
Move-Cell "H1" "G1"
Move-Cell "G1" "F1"
$ws.Range("G1").Value = "{{title}}"
$ws.Range("F1").Clear()

# ---- Row 3: B3:J3 -> C3:K3 (header row) ----
Move-Cell "K3" "J3"
Move-Cell "J3" "I3"
Move-Cell "I3" "H3"
Move-Cell "H3" "G3"
Move-Cell "G3" "F3"
Move-Cell "F3" "E3"
Move-Cell "E3" "D3"
Move-Cell "D3" "C3"
Move-Cell "C3" "B3"

$ws.Range("K3").Value = 3
$ws.Range("J3").Value = 2
$ws.Range("I3").Value = 1
$ws.Range("H3").Value = "Total"
$ws.Range("G3").Value = "{{item}}"
$ws.Range("F3").Value = "Citi"
$ws.Range("E3").Value = "Age"
$ws.Range("D3").Value = "Role"
$ws.Range("C3").Value = "Name"
$ws.Range("B3").Value = "#"

# ---- Row 4: B4:J4 -> C4:K4 (template data row) ----
Move-Cell "K4" "J4"
Move-Cell "J4" "I4"
Move-Cell "I4" "H4"
Move-Cell "H4" "G4"
Move-Cell "G4" "F4"
Move-Cell "F4" "E4"
Move-Cell "E4" "D4"
Move-Cell "D4" "C4"
Move-Cell "C4" "B4"

$ws.Range("K4").Value = 4
$ws.Range("J4").Value = 3
$ws.Range("I4").Value = 2
$ws.Range("H4").Value = "&=SUM(PlanData_Hours)<<OnlyValues>>"
$ws.Range("G4").Value = "{{item}}"
$ws.Range("F4").Value = "{{item.Address.City}}"
$ws.Range("E4").Value = "{{item.Age}}"
$ws.Range("D4").Value = "{{item.Role}}"
$ws.Range("C4").Value = "{{item.Name}}"
$ws.Range("B4").Value = "{{index+1}}"

# ---- Row 5: D5:E5 -> E5:F5 (empty styled cells), H5:J5 -> I5:K5 ----
Move-Cell "K5" "J5"
Move-Cell "J5" "I5"
Move-Cell "I5" "H5"
Move-Cell "F5" "E5"
Move-Cell "E5" "D5"

$ws.Range("K5").Value = 5
$ws.Range("J5").Value = 4
$ws.Range("I5").Value = 3
$ws.Range("D5").Clear()
$ws.Range("H5").Clear()

# ---- Rows 6-9: only the H:J helper strip shifts to I:K; A:G grid is untouched ----
Move-Cell "K6" "J6"
Move-Cell "J6" "I6"
Move-Cell "I6" "H6"
$ws.Range("K6").Value = 6
$ws.Range("J6").Value = 5
$ws.Range("I6").Value = 4
$ws.Range("H6").Clear()

Move-Cell "K7" "J7"
Move-Cell "J7" "I7"
Move-Cell "I7" "H7"
$ws.Range("K7").Value = 7
$ws.Range("J7").Value = 6
$ws.Range("I7").Value = 5
$ws.Range("H7").Clear()

Move-Cell "K8" "J8"
Move-Cell "J8" "I8"
Move-Cell "I8" "H8"
$ws.Range("K8").Value = 8
$ws.Range("J8").Value = 7
$ws.Range("I8").Value = 6
$ws.Range("H8").Clear()

Move-Cell "K9" "J9"
Move-Cell "J9" "I9"
Move-Cell "I9" "H9"
$ws.Range("K9").Value = 9
$ws.Range("J9").Value = 8
$ws.Range("I9").Value = 7
$ws.Range("H9").Clear()

# ---- Sheet-level metadata ----
$ws.Columns("B:H").ColumnWidth = 11.33203125

$ws.Range("B4").Select()

$wb.Names.Add("dates", "=Лист1!`$G`$3")
$wb.Names.Add("PlanData", "=Лист1!`$A`$4:`$H`$5")
$wb.Names.Add("PlanData_Hours", "=Лист1!`$G`$4")
